$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (item #5): new release entry ---

# C10 ("2020.09.18") looks like a date to Excel's normal value parser, which
# would convert it into a serial date number. The source file keeps it as
# literal text (same as the other "Date" column entries), so build it via a
# formula-that-returns-text in a scratch cell and paste back as a value --
# that keeps the literal string (stored as a shared string) without
# disturbing the cell's existing (unformatted) style.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="2020.09.18"'
$scratch.Copy()
$ws.Range("C10").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Range("E10").Value = "20200918_V1_9_Digen_V174_70_FLEET_No_Trigger.7z"
$ws.Range("F10").Value = "임승한"

# G10 ("change" column) wraps text like the rows above it.
$ws.Range("G10").WrapText = $true
$ws.Range("G10").Value = "Debug protocol 9/10 반영"

# Move the active selection to G13, matching the author's final cursor position.
$null = $ws.Range("G13").Select()
